$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.234.29"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "3.749.18"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "594.09"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "169.06"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "3.749.84"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("E8").Value = "  +0.00%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.523"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  +0.29%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.46"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.452"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000274"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.92%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "36.32"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "4.383.35"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "3.760.87"
$ws.Range("E16").Value = "  -2.18%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "18.50"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "67.248.95"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("E20").Value = "  +0.78%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.47"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -5.80%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "467.05"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.715"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -3.04%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "83.71"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.0000147"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -8.72%  "
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "3.904.44"
$ws.Range("E31").Value = "  -2.10%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.63"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.93%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "30.39"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("E34").Value = "  -4.06%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "9.07"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("D36").Value = "3.715.95"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("E39").Value = "  -2.01%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.88%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.81"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -1.90%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.68"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.63%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.93"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("E47").Value = "  -2.78%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "396.08"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -5.51%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.000269"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -8.16%  "
$ws.Range("E50").Value = "  -2.16%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "38.90"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
